$d = $word.ActiveDocument

# --- Paragraph 1: "git add index.html ... /tepuy.css" -> "git status" ---
$p1 = $d.Paragraphs.Item(1).Range
$r1 = $d.Range($p1.Start, $p1.End - 1)
$r1.Text = "git status"

# --- Paragraph 2: git commit -m "Finalised calculators with charts and UI polish" -> git add -A ---
$p2 = $d.Paragraphs.Item(2).Range
$r2 = $d.Range($p2.Start, $p2.End - 1)
$r2.Text = "git add -A"

# --- Paragraph 3: git push -> git commit -m "Deploy full site refresh: ..." ---
$p3 = $d.Paragraphs.Item(3).Range
$r3 = $d.Range($p3.Start, $p3.End - 1)
$r3.Text = 'git commit -m "Deploy full site refresh: improved navbar, mobile menu, styles, images"'

# --- Remove the trailing "git add blog/index.html (fix)" / "git commit ... Fix blog post
#     content and formatting" / "git push" paragraphs (originally paragraphs 15-17) ---
$start = $d.Paragraphs.Item(15).Range.Start
$end = $d.Paragraphs.Item(17).Range.End
$d.Range($start, $end).Delete()

# --- Remove the block between the old "git push" (para 3) and the lone blank paragraph
#     that sits right before the "git add blog/index.html" fix block (originally
#     paragraphs 4-13: four blanks, the "Buy costs" table line, two more blanks, the
#     "git add blog/index.html img/... blog/property-vs-shares-australia.html" line,
#     its commit message, and its "git push") ---
$start2 = $d.Paragraphs.Item(4).Range.Start
$end2 = $d.Paragraphs.Item(13).Range.End
$d.Range($start2, $end2).Delete()

# --- Insert a new paragraph right after the (now third) paragraph for "git push origin main" ---
$afterP3 = $d.Paragraphs.Item(3).Range
$afterP3.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4).Range
$r4 = $d.Range($p4.Start, $p4.End - 1)
$r4.Text = "git push origin main"
